# Update Data Set 5 april 2020
# Append the 5-Apr-2020 (serial 43926) row to each of the three data
# sheets (Confirmed, Recoverd, Death), matching the existing table's
# layout/formatting, then restore the view's selection state.

$wb = $excel.ActiveWorkbook

# ---- Confirmed (sheet1): new row 30 -> 43926, 88, 18 ----
$ws = $wb.Worksheets.Item("Confirmed")
$ws.Range("A29:C29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").Value = 43926
$ws.Range("B30").Value = 88
$ws.Range("C30").Value = 18

# ---- Death (sheet3): new row 30 -> 43926, 9, 1 ----
$ws = $wb.Worksheets.Item("Death")
$ws.Range("A29:C29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").Value = 43926
$ws.Range("B30").Value = 9
$ws.Range("C30").Value = 1

# ---- Recoverd (sheet2): new row 30 -> 43926, 30, 0 ----
$ws = $wb.Worksheets.Item("Recoverd")
$ws.Range("A29:C29").Copy() | Out-Null
$ws.Range("A30").PasteSpecial(-4122) | Out-Null
$ws.Range("A30").Value = 43926
$ws.Range("B30").Value = 30
$ws.Range("C30").Value = 0

# ---- Restore per-sheet selections (activeCell/sqref) ----
$wsConfirmed = $wb.Worksheets.Item("Confirmed")
$wsConfirmed.Select() | Out-Null
$wsConfirmed.Range("B31").Select() | Out-Null

$wsDeath = $wb.Worksheets.Item("Death")
$wsDeath.Select() | Out-Null
$wsDeath.Range("D30").Select() | Out-Null

# Recoverd is the tab that stays active/selected in the workbook, so
# select it last.
$wsRecoverd = $wb.Worksheets.Item("Recoverd")
$wsRecoverd.Select() | Out-Null
$wsRecoverd.Range("D30").Select() | Out-Null
